# Weekly fruit/vegetable data update: insert a new "latest" data row at the
# top of the data block (row 12), pushing every following record down by
# one row (old row 12 -> 13, ..., old row 33 -> 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12; this shifts rows
# 12..33 down to 13..34, preserving their contents and formatting.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 44519
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 300000000
$ws.Range("G12").Value = "Espárragos"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1600
$ws.Range("L12").Value = 1800
$ws.Range("M12").Value = 1700
$ws.Range("N12").Value = "`$/kilo"
$ws.Range("O12").Value = "Provincia de Linares"
$ws.Range("P12").Value = 1700
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
